$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the progress value of row 24 (alinea 1p) to 100 -> marks it as "Done!"
$ws.Range("D24").Value = 100
